$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 97: change B97 to a numeric value (3) instead of text "3"
$ws.Range("B97").Value = 3

# Insert new row 98 with data
$ws.Range("A98").Value = "Ying Tang"
$ws.Range("B98").NumberFormat = "@"
$ws.Range("B98").Value = "3"
$ws.Range("C98").Value = "无"
$ws.Range("D98").Value = "ACK"
$ws.Range("E98").Value = "EXP"
$ws.Range("F98").Value = "840f898f-6d0b-4603-abaa-7e0871215f61"
$ws.Range("G98").Value = "HyEi7bWR-_annotated.xlsx"
$ws.Range("H98").Value = "We have carried out additional experiments to examine run time and the following results will be included in the revision."
